# Auto-generated edit script: updates FFXIV Lich Profits market-data cells (H:N)
# across sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR to match the scheduled
# market-data refresh captured in the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 1097.2632
$ws.Range("I39").Value = 944.1
$ws.Range("J39").Value = 1267.4445
$ws.Range("K39").Value = 2832.3
$ws.Range("L39").Value = 3802.3335
$ws.Range("M39").Value = -2536.3
$ws.Range("N39").Value = -4394.333500000001
$ws.Range("H70").Value = 11858.97
$ws.Range("I70").Value = 3015.8333
$ws.Range("J70").Value = 13824.111
$ws.Range("K70").Value = 9047.499899999999
$ws.Range("L70").Value = 41472.333
$ws.Range("M70").Value = -8777.499899999999
$ws.Range("N70").Value = -42012.333
$ws.Range("H73").Value = 11858.97
$ws.Range("I73").Value = 3015.8333
$ws.Range("J73").Value = 13824.111
$ws.Range("K73").Value = 9047.499899999999
$ws.Range("L73").Value = 41472.333
$ws.Range("M73").Value = -8111.499899999999
$ws.Range("N73").Value = -43344.333
$ws.Range("H106").Value = 101675.5
$ws.Range("I106").Value = 101675.5
$ws.Range("K106").Value = 101675.5
$ws.Range("M106").Value = -101044.5
$ws.Range("H112").Value = 3074985.2
$ws.Range("J112").Value = 3928917.8
$ws.Range("L112").Value = 11786753.4
$ws.Range("N112").Value = -11788969.4
$ws.Range("H138").Value = 2029.6061
$ws.Range("I138").Value = 1735.8462
$ws.Range("J138").Value = 2220.55
$ws.Range("K138").Value = 5207.5386
$ws.Range("L138").Value = 6661.650000000001
$ws.Range("M138").Value = -67.53859999999986
$ws.Range("N138").Value = -16941.65

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17131.66
$ws.Range("I32").Value = 16405.531
$ws.Range("J32").Value = 28507.666
$ws.Range("K32").Value = 16405.531
$ws.Range("L32").Value = 28507.666
$ws.Range("M32").Value = -16118.531
$ws.Range("N32").Value = -29081.666
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("H61").Value = 3033.5588
$ws.Range("I61").Value = 1254
$ws.Range("K61").Value = 1254
$ws.Range("M61").Value = -1042
$ws.Range("H110").Value = 2723.125
$ws.Range("I110").Value = 1489.4
$ws.Range("K110").Value = 1489.4
$ws.Range("M110").Value = 555.5999999999999
$ws.Range("H122").Value = 9318
$ws.Range("I122").Value = 5824.5557
$ws.Range("K122").Value = 17473.6671
$ws.Range("M122").Value = -15023.6671
$ws.Range("H132").Value = 3051.2068
$ws.Range("I132").Value = 2141.2727
$ws.Range("K132").Value = 6423.8181
$ws.Range("M132").Value = -3893.8181
$ws.Range("H136").Value = 3033.5588
$ws.Range("I136").Value = 1254
$ws.Range("K136").Value = 3762
$ws.Range("M136").Value = -1212

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H29").Value = 11158.2
$ws.Range("I29").Value = 900
$ws.Range("J29").Value = 17997
$ws.Range("K29").Value = 900
$ws.Range("L29").Value = 17997
$ws.Range("M29").Value = -611
$ws.Range("N29").Value = -18575
$ws.Range("H54").Value = 6750
$ws.Range("I54").Value = 3500
$ws.Range("J54").Value = 10000
$ws.Range("K54").Value = 3500
$ws.Range("L54").Value = 10000
$ws.Range("M54").Value = -3016
$ws.Range("N54").Value = -10968
$ws.Range("H86").Value = 1689.3572
$ws.Range("I86").Value = 1395.6818
$ws.Range("J86").Value = 2766.1667
$ws.Range("K86").Value = 1395.6818
$ws.Range("L86").Value = 2766.1667
$ws.Range("M86").Value = -272.6818000000001
$ws.Range("N86").Value = -5012.1667
$ws.Range("H89").Value = 1689.3572
$ws.Range("I89").Value = 1395.6818
$ws.Range("J89").Value = 2766.1667
$ws.Range("K89").Value = 6978.409000000001
$ws.Range("L89").Value = 13830.8335
$ws.Range("M89").Value = -1362.409000000001
$ws.Range("N89").Value = -25062.8335
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 2848.125
$ws.Range("I10").Value = 1380.8334
$ws.Range("J10").Value = 7250
$ws.Range("K10").Value = 1380.8334
$ws.Range("L10").Value = 7250
$ws.Range("M10").Value = -1241.8334
$ws.Range("N10").Value = -7528
$ws.Range("H31").Value = 18573.166
$ws.Range("I31").Value = 1436.6451
$ws.Range("J31").Value = 49822.117
$ws.Range("K31").Value = 1436.6451
$ws.Range("L31").Value = 49822.117
$ws.Range("M31").Value = -1141.6451
$ws.Range("N31").Value = -50412.117
$ws.Range("H34").Value = 18573.166
$ws.Range("I34").Value = 1436.6451
$ws.Range("J34").Value = 49822.117
$ws.Range("K34").Value = 1436.6451
$ws.Range("L34").Value = 49822.117
$ws.Range("M34").Value = -1234.6451
$ws.Range("N34").Value = -50226.117
$ws.Range("H39").Value = 38375
$ws.Range("I39").Value = 3500
$ws.Range("K39").Value = 3500
$ws.Range("M39").Value = -3109
$ws.Range("H49").Value = 38375
$ws.Range("I49").Value = 3500
$ws.Range("K49").Value = 3500
$ws.Range("M49").Value = -3318

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 6250135.5
$ws.Range("I2").Value = 60.666668
$ws.Range("J2").Value = 14285946
$ws.Range("K2").Value = 364.000008
$ws.Range("L2").Value = 85715676
$ws.Range("M2").Value = -251.000008
$ws.Range("N2").Value = -85715902
$ws.Range("H17").Value = 113.166664
$ws.Range("J17").Value = 93
$ws.Range("L17").Value = 279
$ws.Range("N17").Value = -617
$ws.Range("H23").Value = 1288.5
$ws.Range("J23").Value = 1471.421
$ws.Range("L23").Value = 4414.263
$ws.Range("N23").Value = -4884.263
$ws.Range("H38").Value = 149.22223
$ws.Range("I38").Value = 175.66667
$ws.Range("J38").Value = 136
$ws.Range("K38").Value = 527.00001
$ws.Range("L38").Value = 408
$ws.Range("M38").Value = -180.00001
$ws.Range("N38").Value = -1102
$ws.Range("H44").Value = 3095.8572
$ws.Range("I44").Value = 500
$ws.Range("K44").Value = 1500
$ws.Range("M44").Value = -1102
$ws.Range("H86").Value = 494.8
$ws.Range("I86").Value = 494.8
$ws.Range("K86").Value = 1484.4
$ws.Range("M86").Value = -298.4000000000001
$ws.Range("H89").Value = 494.8
$ws.Range("I89").Value = 494.8
$ws.Range("K89").Value = 4453.2
$ws.Range("M89").Value = 1474.8
$ws.Range("H131").Value = 11629576
$ws.Range("J131").Value = 2025
$ws.Range("L131").Value = 6075
$ws.Range("N131").Value = -16155

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 83.19047500000001
$ws.Range("I2").Value = 99.15385000000001
$ws.Range("J2").Value = 57.25
$ws.Range("K2").Value = 99.15385000000001
$ws.Range("L2").Value = 57.25
$ws.Range("M2").Value = 13.84614999999999
$ws.Range("N2").Value = -283.25
$ws.Range("H31").Value = 1312.2
$ws.Range("I31").Value = 1312.2
$ws.Range("K31").Value = 1312.2
$ws.Range("M31").Value = -1020.2
$ws.Range("H37").Value = 1312.2
$ws.Range("I37").Value = 1312.2
$ws.Range("K37").Value = 1312.2
$ws.Range("M37").Value = -1035.2
$ws.Range("H122").Value = 1583.3889
$ws.Range("I122").Value = 1426
$ws.Range("K122").Value = 4278
$ws.Range("M122").Value = -1828
$ws.Range("H132").Value = 42974.383
$ws.Range("I132").Value = 45971.043
$ws.Range("K132").Value = 137913.129
$ws.Range("M132").Value = -135383.129
$ws.Range("H136").Value = 27173.867
$ws.Range("J136").Value = 27173.867
$ws.Range("L136").Value = 81521.601
$ws.Range("N136").Value = -86621.601

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7079.6665
$ws.Range("I7").Value = 7163.0605
$ws.Range("J7").Value = 6773.8887
$ws.Range("K7").Value = 7163.0605
$ws.Range("L7").Value = 6773.8887
$ws.Range("M7").Value = -7051.0605
$ws.Range("N7").Value = -6997.8887
$ws.Range("H22").Value = 2091.5789
$ws.Range("I22").Value = 1442.3077
$ws.Range("K22").Value = 1442.3077
$ws.Range("M22").Value = -1147.3077
$ws.Range("H27").Value = 2091.5789
$ws.Range("I27").Value = 1442.3077
$ws.Range("K27").Value = 1442.3077
$ws.Range("M27").Value = -1335.3077
$ws.Range("H36").Value = 88888
$ws.Range("J36").Value = 88888
$ws.Range("L36").Value = 88888
$ws.Range("N36").Value = -90012
$ws.Range("H41").Value = 43343.668
$ws.Range("I41").Value = 43343.668
$ws.Range("K41").Value = 43343.668
$ws.Range("M41").Value = -42905.668
$ws.Range("H46").Value = 3319.8823
$ws.Range("I46").Value = 1387.6154
$ws.Range("K46").Value = 1387.6154
$ws.Range("M46").Value = -1199.6154
$ws.Range("H122").Value = 4160.8335
$ws.Range("I122").Value = 3813.889
$ws.Range("K122").Value = 11441.667
$ws.Range("M122").Value = -8991.667000000001
$ws.Range("H126").Value = 7079.6665
$ws.Range("I126").Value = 7163.0605
$ws.Range("J126").Value = 6773.8887
$ws.Range("K126").Value = 21489.1815
$ws.Range("L126").Value = 20321.6661
$ws.Range("M126").Value = -19019.1815
$ws.Range("N126").Value = -25261.6661
$ws.Range("H136").Value = 2315.1162
$ws.Range("I136").Value = 1852.0571
$ws.Range("K136").Value = 5556.1713
$ws.Range("M136").Value = -3006.1713

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H42").Value = 44832.668
$ws.Range("I42").Value = 44832.668
$ws.Range("K42").Value = 44832.668
$ws.Range("M42").Value = -44454.668
$ws.Range("H70").Value = 39899.5
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()
$ws.Range("H73").Value = 39899.5
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()
$ws.Range("H135").Value = 54332.777
$ws.Range("J135").Value = 56124.375
$ws.Range("L135").Value = 56124.375
$ws.Range("N135").Value = -66264.375
